$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 13 (old extra fixtures), keeping only rows 1-2
$ws.Range("A3:B13").Delete()

# Update row 1
$ws.Range("A1").Value = "Manchester United Legends v Liverpool Legends "
$ws.Range("B1").Value = "21 MaySat15:00"

# Update row 2
$ws.Range("A2").Value = "Women's EURO 2022 England v Austria "
$ws.Range("B2").Value = "06 JulWed20:00"
